$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.602.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.925.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4057"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08225"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.954.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.088"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.255"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06877"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.012"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.580.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.189"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.171.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.437"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.012"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09631"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.631"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.581"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06413"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.185"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5953"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.884"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1852"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.435"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.255"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07531"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5559"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.950"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.436"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.26%  "
